# Update the cryptos list with the latest scraped data.
# Mirrors the GitHub Actions scheduled job that refreshes prices / 1h volume
# and occasionally re-sorts two rows when their ranking swaps places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D prices are plain text (e.g. "67.361.85" thousand-grouped, or
# "597.07" decimal). Values that look like a genuine number would otherwise
# get auto-converted by Excel, so force those specific cells to Text first
# so the write lands as a literal string, matching the source feed.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.361.85"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "3.528.97"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "597.07"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "173.71"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").Value = "0.133"
$ws.Range("E9").Value = "  +7.42%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "4.138.31"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "0.0000181"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "67.273.61"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "3.534.80"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "14.18"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").Value = "396.67"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "73.63"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.540"
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "10.23"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "6.30"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "2.08"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").Value = "24.08"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("D35").Value = "163.98"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "6.87"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "4.73"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").Value = "0.0747"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "26.54"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "27.23"
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  +4.82%  "
$ws.Range("D44").Value = "2.816.07"
$ws.Range("D45").Value = "42.98"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").Value = "0.0310"
$ws.Range("D47").Value = "343.00"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "33.90"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").Value = "6.53"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "0.853"
$ws.Range("E51").Value = "  +0.53%  "
